# CBPX_QTR_FIN.xlsx - "Doing Updates for Financials"
# Insert two new quarterly-result columns (D:E) in front of the existing
# D:K data block on the CBPX sheet, shifting the old D:K data right to
# F:M, and populate the two new columns with the newest two quarters of
# figures (plus a small number of restated prior-quarter figures in the
# "Capital Expenditures" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before column D; this shifts the existing
#    D:K data block to F:M and keeps formulas / column widths in sync.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) The freshly-inserted D:E columns have no number formatting yet.
#    Copy the formatting (date format row 7/38/80, plain-number format
#    everywhere else, etc.) from the new column F (the old column D,
#    now shifted one column over) so D:E look exactly like the other
#    quarter columns.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Write the two newest quarters of data into the new D (most recent)
#    and E (second most recent) columns for every data row.
$newData = @(
  @(7, 43465, 43373),
  @(8, 140800, 131200),
  @(9, 101200, 94300),
  @(10, 39600, 36900),
  @(12, "NA", "NA"),
  @(13, 0, 0),
  @(14, 0, 0),
  @(15, 0, 0),
  @(17, 111700, 104300),
  @(18, 29100, 26900),
  @(20, -3100, -2900),
  @(21, 36700, 35600),
  @(22, 0, 0),
  @(23, 25900, 24000),
  @(24, 5800, 5400),
  @(25, 0, 0),
  @(26, 20100, 18600),
  @(27, 20100, 18600),
  @(28, 0, 0),
  @(29, 0, "NA"),
  @(30, 0, 0),
  @(31, 0, 0),
  @(32, 3100, 2900),
  @(33, 20100, 18600),
  @(34, 0, 0),
  @(35, 20100, 18600),
  @(38, 43465, 43373),
  @(41, 102600, 105500),
  @(42, 0, 0),
  @(43, 38500, 39600),
  @(44, 32200, 32400),
  @(45, 19800, 11600),
  @(46, 193100, 189100),
  @(47, 8000, 8200),
  @(48, 288400, 290700),
  @(49, 182600, 184600),
  @(50, 0, 0),
  @(51, 0, 0),
  @(52, 300, 300),
  @(53, 0, 0),
  @(54, 672400, 672900),
  @(57, 48100, 32300),
  @(58, 1700, 1700),
  @(59, 12800, 13200),
  @(60, 62500, 47200),
  @(61, 261900, 262400),
  @(62, 20200, 15400),
  @(63, 0, 0),
  @(64, 0, 0),
  @(65, 0, 0),
  @(66, 344600, 325000),
  @(68, 0, 0),
  @(69, 0, 0),
  @(70, 0, 0),
  @(71, 0, 0),
  @(72, 212600, 192700),
  @(73, 0, 0),
  @(74, 0, 0),
  @(75, 0, 0),
  @(76, 327700, 348000),
  @(77, 0, 0),
  @(80, 43465, 43373),
  @(81, 20100, 18600),
  @(83, 10800, 11600),
  @(84, 0, 0),
  @(85, 0, 0),
  @(86, 0, 0),
  @(87, 0, 0),
  @(88, 0, 0),
  @(89, 46800, 31300),
  @(91, -9100, -6800),
  @(92, 0, 0),
  @(93, 0, 0),
  @(94, -9800, -6900),
  @(96, 0, 0),
  @(97, 0, 0),
  @(98, 0, 0),
  @(99, 0, 0),
  @(100, -39400, -3400),
  @(101, -400, 100),
  @(102, -2800, 21100)
)

foreach ($item in $newData) {
    $r = $item[0]
    $ws.Cells.Item($r, 4).Value2 = $item[1]   # column D
    $ws.Cells.Item($r, 5).Value2 = $item[2]   # column E
}

# 4) Row 91 ("Capital Expenditures") also had several of the already
#    shifted prior-quarter values restated; overwrite F:K with the
#    corrected figures (L:M keep the values carried over from the shift).
$ws.Cells.Item(91, 6).Value2 = -7100   # F91
$ws.Cells.Item(91, 7).Value2 = -6000   # G91
$ws.Cells.Item(91, 8).Value2 = -7400   # H91
$ws.Cells.Item(91, 9).Value2 = -6000   # I91
$ws.Cells.Item(91, 10).Value2 = -2700  # J91
$ws.Cells.Item(91, 11).Value2 = -5400  # K91
